$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.830.76"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "3.058.20"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.97"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.26"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.056.59"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.10"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.480"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.28"
$ws.Range("D15").Value = "3.559.31"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "63.866.95"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "3.062.24"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "487.17"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.33"
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.70"
$ws.Range("E23").Value = "  +8.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.52"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.62"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.46"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("E32").Value = "  +5.44%  "
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.04"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "441.66"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0814"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "3.043.29"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -7.01%  "
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.276"
$ws.Range("E44").Value = "  +7.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.08"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  +6.41%  "
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.49"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +4.52%  "
